$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B/C/E (and any non-numeric-looking D) text updates ---
$ws.Range("D2").Value = '28.124.53'
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").Value = '1.895.38'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5023'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3897'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09226'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.129'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.391'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.38%  '
$ws.Range("E13").Value = '  -1.74%  '
$ws.Range("D14").Value = '1.903.95'
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.283'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.77%  '
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("E17").Value = '  -2.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06652'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.22%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.207'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.45%  '
$ws.Range("D23").Value = '28.186.47'
$ws.Range("E23").Value = '  -1.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  +1.86%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.375'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = '2.124.59'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.553'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.63%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '158.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.03%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '126.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.57%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.082'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.24%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.1057'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.608'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.613'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.557'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.355'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.84%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06599'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.98%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02401'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.85%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2205'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.224'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.71%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6472'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("B44").Value = 'InternetComputer(DFINITY)'
$ws.Range("C44").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.971'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.83%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6105'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '13.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.17%  '
$ws.Range("B48").Value = 'WEMIXTOKEN'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.306'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.85%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.693'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.001'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.17%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '122.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.14%  '
